$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 427, shifting existing rows 427:508 down to 428:509
$ws.Rows.Item(427).Insert()

# Populate the newly inserted row 427 with a new weekly price record,
# matching the constant columns used throughout this data block.
$ws.Cells.Item(427, 1).Value = 3
$ws.Cells.Item(427, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(427, 3).Value = "Coquimbo"
$ws.Cells.Item(427, 4).Value = 44995
$ws.Cells.Item(427, 5).Value = 5
$ws.Cells.Item(427, 6).Value = 100112009
$ws.Cells.Item(427, 7).Value = "Acelga"
$ws.Cells.Item(427, 8).Value = "Sin especificar"
$ws.Cells.Item(427, 9).Value = "Primera"
$ws.Cells.Item(427, 10).Value = 230
$ws.Cells.Item(427, 11).Value = 3300
$ws.Cells.Item(427, 12).Value = 3500
$ws.Cells.Item(427, 13).Value = 3404
$ws.Cells.Item(427, 14).Value = "`$/docena de atados (6 kilos)"
$ws.Cells.Item(427, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(427, 16).Value = 567
$ws.Cells.Item(427, 17).Value = 6
$ws.Cells.Item(427, 18).Value = "Hortaliza"
